# Price update for 2026-02-07
# Append a new tracked data-point row to the bottom of the price-history
# sheet (Date, Price, Discount, Incredible).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$dataRange = $ws.Range("A" + $newRow + ":D" + $newRow)

# Force the new cells to be stored as plain text (matching how every other
# history row is stored, i.e. shared strings) instead of letting Excel
# auto-detect the date-looking / numeric-looking text and silently convert
# it into a real date serial or number.
$dataRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "647000"
$ws.Cells.Item($newRow, 3).Value = "28"
$ws.Cells.Item($newRow, 4).Value = "0"

# Restore the default "Normal" style so the new row doesn't carry any
# extra formatting beyond what the rest of the sheet already uses.
$dataRange.Style = "Normal"
